$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Targeted single-cell corrections ---
$ws.Cells.Item(60, 17).Value = 0      # Q60: 1 -> 0
$ws.Cells.Item(1056, 15).Value = 2    # O1056: 0 -> 2
$ws.Cells.Item(1058, 18).Value = 0    # R1058: blank -> 0
$ws.Cells.Item(1059, 18).Value = 0    # R1059: blank -> 0

# --- New weekly rows 1060-1072 (appended history) ---
$newRows = @(
    @{ Row=1060; Vals=@(45474, 333, 341, 330.25, 337.3999938964844, 334.4373168945312, 18454292, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0) },
    @{ Row=1061; Vals=@(45481, 337.0499877929688, 343, 321.2000122070312, 339.0499877929688, 336.0728149414062, 19977100, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0) },
    @{ Row=1062; Vals=@(45488, 341.4500122070312, 355.75, 338.25, 342.8999938964844, 342.8999938964844, 18606630, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1) },
    @{ Row=1063; Vals=@(45495, 342.2999877929688, 377.9500122070312, 329.8999938964844, 372.3999938964844, 372.3999938964844, 47792130, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0) },
    @{ Row=1064; Vals=@(45502, 370.9500122070312, 376.4500122070312, 360.25, 362.3999938964844, 362.3999938964844, 22435071, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0) },
    @{ Row=1065; Vals=@(45509, 354, 370.75, 350.5, 366.7000122070312, 366.7000122070312, 20923063, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0) },
    @{ Row=1066; Vals=@(45516, 364.3500061035156, 375.6000061035156, 361.8500061035156, 370.8999938964844, 370.8999938964844, 11338898, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0) },
    @{ Row=1067; Vals=@(45523, 370, 384.2000122070312, 370, 370.7999877929688, 370.7999877929688, 15726586, 2024, 8, 19, 0, 0, 0, 34, 1, 0, 0) },
    @{ Row=1068; Vals=@(45530, 372.7999877929688, 372.7999877929688, 361.0499877929688, 367.2999877929688, 367.2999877929688, 13438876, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0) },
    @{ Row=1069; Vals=@(45537, 369, 372.25, 352, 354, 354, 11326345, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0) },
    @{ Row=1070; Vals=@(45544, 348.2000122070312, 352.3999938964844, 334.1000061035156, 334.8500061035156, 334.8500061035156, 14646069, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0) },
    @{ Row=1071; Vals=@(45551, 335, 338.2999877929688, 319.7999877929688, 325, 325, 20467126, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0) },
    @{ Row=1072; Vals=@(45558, 327, 342.3999938964844, 325.5499877929688, 339.3500061035156, 339.3500061035156, 21589457, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0) }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $vals = $r.Vals
    # Column A (Datetime serial) with the date-style used by the rest of the column
    $aCell = $ws.Cells.Item($rowNum, 1)
    $aCell.Value = $vals[0]
    $aCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    # Columns B..Q (index 1..16 in $vals maps to spreadsheet columns 2..17)
    for ($c = 1; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $vals[$c]
    }
    # Column R (backup) stays blank, matching the source data
    $ws.Cells.Item($rowNum, 18).Value = ""
}

Write-Output "done"
